$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "30.739.71"
$ws.Range("E2").Value = "  +2.56%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.894.55"
$ws.Range("E3").Value = "  +0.91%  "

# Row 4 - TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  +0.20%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.03"
$ws.Range("E5").Value = "  +1.77%  "

# Row 6 - USDC
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("E6").Value = "  +0.21%  "

# Row 7 - XRP
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4933"
$ws.Range("E7").Value = "  -1.08%  "

# Row 8 - Cardano
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2958"
$ws.Range("E8").Value = "  +1.26%  "

# Row 9 - Dogecoin
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06821"
$ws.Range("E9").Value = "  +3.18%  "

# Row 10 - WrappedEther
$ws.Range("D10").Value = "1.895.77"
$ws.Range("E10").Value = "  +0.95%  "

# Row 11 - Solana
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "17.39"
$ws.Range("E11").Value = "  +4.01%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  +0.11%  "

# Row 13 - Litecoin
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "92.11"
$ws.Range("E13").Value = "  +7.01%  "

# Row 14 - Polkadot -> Polygon
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6871"
$ws.Range("E14").Value = "  +3.06%  "

# Row 15 - Polygon -> Polkadot
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.105"
$ws.Range("E15").Value = "  +5.14%  "

# Row 16 - WrappedBTC
$ws.Range("D16").Value = "30.713.90"
$ws.Range("E16").Value = "  +2.59%  "

# Row 17 - ShibaInu
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008004"
$ws.Range("E17").Value = "  +1.33%  "

# Row 18 - Avalanche
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.39"
$ws.Range("E18").Value = "  +5.06%  "

# Row 19 - Dai
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9995"
$ws.Range("E19").Value = "  +0.13%  "

# Row 20 - WrappedliquidstakedEther2.0
$ws.Range("D20").Value = "2.138.81"
$ws.Range("E20").Value = "  +0.94%  "

# Row 21 - BinanceUSD
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  +0.36%  "

# Row 22 - Uniswap
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.869"
$ws.Range("E22").Value = "  +2.44%  "

# Row 23 - BitcoinCash
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "191.03"
$ws.Range("E23").Value = "  +36.10%  "

# Row 24 - Chainlink
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.095"
$ws.Range("E24").Value = "  +8.10%  "

# Row 25 - Cosmos
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.403"
$ws.Range("E25").Value = "  +3.90%  "

# Row 26 - Monero
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.37"
$ws.Range("E26").Value = "  +4.29%  "

# Row 27 - EthereumClassic
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.33"
$ws.Range("E27").Value = "  +13.71%  "

# Row 28 - LidoDAOToken
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.935"
$ws.Range("E28").Value = "  +1.41%  "

# Row 29 - Toncoin
$ws.Range("E29").Value = "  +0.56%  "

# Row 30 - InternetComputer(DFINITY)
$ws.Range("E30").Value = "  +4.84%  "

# Row 31 - Stellar
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09031"
$ws.Range("E31").Value = "  +2.83%  "

# Row 32 - Filecoin
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.057"
$ws.Range("E32").Value = "  +2.70%  "

# Row 33 - Hedera
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05198"
$ws.Range("E33").Value = "  +2.70%  "

# Row 34 - ImmutableX
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7501"
$ws.Range("E34").Value = "  +5.29%  "

# Row 35 - ARBITRUM
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.131"
$ws.Range("E35").Value = "  +2.21%  "

# Row 36 - HuobiToken
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.709"
$ws.Range("E36").Value = "  +1.69%  "

# Row 37 - VeChain
$ws.Range("E37").Value = "  +6.52%  "

# Row 38 - MXToken
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.676"
$ws.Range("E38").Value = "  -0.36%  "

# Row 39 - RenderToken
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.174"
$ws.Range("E39").Value = "  -0.35%  "

# Row 40 - TrustWalletToken
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9392"
$ws.Range("E40").Value = "  +1.05%  "

# Row 41 - TheSandbox
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4454"
$ws.Range("E41").Value = "  +4.76%  "

# Row 42 - Quant
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "106.46"
$ws.Range("E42").Value = "  +4.50%  "

# Row 43 - FraxShare
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.839"
$ws.Range("E43").Value = "  +0.89%  "

# Row 44 - PaxDollar
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9998"
$ws.Range("E44").Value = "  +0.27%  "

# Row 45 - Aptos
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.738"
$ws.Range("E45").Value = "  +3.92%  "

# Row 46 - Algorand
$ws.Range("E46").Value = "  +7.30%  "

# Row 47 - Cronos
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05853"
$ws.Range("E47").Value = "  +3.52%  "

# Row 48 - EnergySwap
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.830"
$ws.Range("E48").Value = "  +7.42%  "

# Row 49 - Decentraland
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.3977"
$ws.Range("E49").Value = "  +5.97%  "

# Row 50 - NEARProtocol
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.428"
$ws.Range("E50").Value = "  +7.32%  "

# Row 51 - Elrond
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "33.65"
$ws.Range("E51").Value = "  +3.88%  "
